$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1012
$ws.Range("F3").Value = 576
$ws.Range("F4").Value = 9056
$ws.Range("F5").Value = 185
$ws.Range("F6").Value = 58
$ws.Range("F7").Value = 1952
$ws.Range("F8").Value = 6285
$ws.Range("F9").Value = 612
$ws.Range("F12").Value = 9313
$ws.Range("F13").Value = 10811
$ws.Range("F14").Value = 1214
$ws.Range("G14").Value = 168
$ws.Range("F15").Value = 1100
$ws.Range("F16").Value = 4863
$ws.Range("F17").Value = 777
$ws.Range("F18").Value = 424
$ws.Range("F22").Value = 1315
$ws.Range("F23").Value = 222
$ws.Range("F24").Value = 1847
$ws.Range("F25").Value = 860
$ws.Range("F26").Value = 1192
$ws.Range("F27").Value = 848
$ws.Range("F28").Value = 2004
$ws.Range("F29").Value = 408
$ws.Range("F30").Value = 591
$ws.Range("F31").Value = 2595
$ws.Range("F33").Value = 173
$ws.Range("F34").Value = 1693
$ws.Range("F35").Value = 90
$ws.Range("F36").Value = 1326
$ws.Range("F37").Value = 421
$ws.Range("F38").Value = 12
$ws.Range("F39").Value = 899
$ws.Range("F40").Value = 575
$ws.Range("F41").Value = 3264
$ws.Range("F42").Value = 231
$ws.Range("F44").Value = 493
$ws.Range("F45").Value = 563
$ws.Range("F47").Value = 890

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5781

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1012
$ws.Range("F3").Value = 576
$ws.Range("F4").Value = 9056
$ws.Range("F5").Value = 58
$ws.Range("F8").Value = 6285
$ws.Range("F9").Value = 612
$ws.Range("F10").Value = 9314
$ws.Range("F11").Value = 9314
$ws.Range("F12").Value = 10811
$ws.Range("F14").Value = 1214
$ws.Range("G14").Value = 168
$ws.Range("F15").Value = 1100
$ws.Range("F16").Value = 4863
$ws.Range("F17").Value = 777
$ws.Range("F18").Value = 424
$ws.Range("F22").Value = 1315
$ws.Range("F23").Value = 222
$ws.Range("F24").Value = 860
$ws.Range("F25").Value = 1192
$ws.Range("F26").Value = 848
$ws.Range("F28").Value = 2004
$ws.Range("F29").Value = 409
$ws.Range("F30").Value = 2595
$ws.Range("F31").Value = 173
$ws.Range("F32").Value = 1693
$ws.Range("F33").Value = 90
$ws.Range("F35").Value = 421
$ws.Range("F39").Value = 899
$ws.Range("F40").Value = 575
$ws.Range("F42").Value = 231
$ws.Range("F44").Value = 493
$ws.Range("F45").Value = 563
$ws.Range("F46").Value = 890

